$wb = $excel.ActiveWorkbook

# --- Add the four new "backward pass" sheets by duplicating the matching
# --- forward-pass sheets (same shape/style/values), then renaming them. ---

$srcNames = @(
    "get_fwd_release_delays",
    "get_fwd_proc_compute_node",
    "get_fwd_end_local",
    "get_trans_back"
)
$newNames = @(
    "get_bwd_release_delays",
    "get_bwd_proc_compute_node",
    "get_bwd_end_local",
    "get_grad_trans_back"
)

for ($i = 0; $i -lt $srcNames.Length; $i++) {
    $src = $wb.Worksheets.Item($srcNames[$i])
    $afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $src.Copy([System.Reflection.Missing]::Value, $afterSheet)
    $created = $wb.Worksheets.Item($wb.Worksheets.Count)
    $created.Name = $newNames[$i]
}

# --- Fix up the selection on each new sheet (copy carries over the
# --- source sheet's selection, so every new sheet needs its own). ---

$wb.Worksheets.Item("get_bwd_release_delays").Range("Q35").Select()
$wb.Worksheets.Item("get_bwd_proc_compute_node").Range("A3").Select()
$wb.Worksheets.Item("get_bwd_end_local").Range("O30").Select()
$wb.Worksheets.Item("get_grad_trans_back").Range("B8").Select()

# --- The active/selected tab ends up on the new "get_bwd_end_local" sheet. ---

$wb.Worksheets.Item("get_bwd_end_local").Activate()
